# Connor Readnour resume — content edits
#
# Two bullet points in the "Experience" section were re-worded:
#   1. "Designed and executed thousands of promotional ..." ->
#      "Designed and executed over 10,000 promotional ..." (and the
#      "utm parameters, and fallbacks" -> "UTM parameters and fallbacks"
#      wording/casing/comma tidy-up).
#   2. The "Demonstrated strong communication skills ..." bullet was
#      re-typed (same final wording) which, in the saved file, shows up
#      as the three runs that made up that sentence being re-combined
#      into a single run.
#
# We use Find to *locate* the text and then assign straight to
# Range.Text (rather than passing the replacement into Find.Execute's
# Replace argument) so that Word's "smart quotes" AutoCorrect doesn't
# mangle the straight apostrophe in "requestor's" into a curly one —
# Range.Text is a literal assignment with no AutoCorrect post-processing.

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $newText
        return $true
    }
    return $false
}

# 1) "... thousands of promotional ... utm parameters, and fallbacks ..."
$old1 = "Designed and executed thousands of promotional Email, Journey, and Automation test cases, ensuring seamless functionality of all links, CTAs, dynamic data, utm parameters, and fallbacks within each email, in accordance with the requestor's specifications."
$new1 = "Designed and executed over 10,000 promotional Email, Journey, and Automation test cases, ensuring seamless functionality of all links, CTAs, dynamic data, UTM parameters and fallbacks within each email, in accordance with the requestor's specifications."
$ok1 = Replace-ExactText $old1 $new1
Write-Output ("bullet 1 updated: " + $ok1)

# 2) "Demonstrated strong communication skills ..." — re-write with the
# identical final text so the three runs collapse back into one. This
# sentence has no apostrophes/quotes, so it's safe to let Find.Execute
# perform the replacement directly (wdReplaceAll) — doing so rebuilds a
# single run for the whole match instead of leaving the original three
# runs in place.
$old2 = "Demonstrated strong communication skills by consistently delivering weekly QA Status Update emails, presenting comprehensive test case status reports from all team members, and effectively highlighting any potential blockers faced by the QA team."
$new2 = "Demonstrated strong communication skills by consistently delivering weekly QA Status Update emails, presenting comprehensive test case status reports from all team members, and effectively highlighting any potential blockers faced by the QA team."
$ok2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Output ("bullet 2 normalized: " + $ok2)
